{"js": "const replacements = [\n  [\"2024-01-22 Monday\", \"2024-01-23 Tuesday\"],\n  [\"604\u00f73=201, 1\", \"908\u00f72=454, 0\"],\n  [\"961\u00f75=192, 1\", \"884\u00f73=294, 2\"],\n  [\"951\u00f79=105, 6\", \"790\u00f78=98, 6\"],\n  [\"944\u00f72=472, 0\", \"620\u00f78=77, 4\"],\n  [\"761\u00f79=84, 5\", \"963\u00f73=321, 0\"],\n  [\"976\u00f78=122, 0\", \"117\u00f77=16, 5\"],\n  [\"309\u00f78=38, 5\", \"729\u00f78=91, 1\"],\n  [\"610\u00f78=76, 2\", \"775\u00f75=155, 0\"],\n  [\"209\u00f77=29, 6\", \"435\u00f74=108, 3\"],\n  [\"831\u00f76=138, 3\", \"933\u00f79=103, 6\"],\n  [\"561\u00f76=93, 3\", \"166\u00f75=33, 1\"],\n  [\"298\u00f74=74, 2\", \"431\u00f73=143, 2\"],\n  [\"359\u00f72=179, 1\", \"409\u00f78=51, 1\"],\n  [\"901\u00f72=450, 1\", \"830\u00f73=276, 2\"],\n  [\"962\u00f73=320, 2\", \"289\u00f74=72, 1\"],\n  [\"142\u00f72=71, 0\", \"443\u00f74=110, 3\"],\n  [\"373\u00f79=41, 4\", \"633\u00f77=90, 3\"],\n  [\"770\u00f74=192, 2\", \"751\u00f75=150, 1\"],\n  [\"326\u00f73=108, 2\", \"325\u00f74=81, 1\"],\n  [\"734\u00f78=91, 6\", \"140\u00f76=23, 2\"],\n  [\"228\u00f78=28, 4\", \"171\u00f77=24, 3\"],\n  [\"810\u00f79=90, 0\", \"488\u00f78=61, 0\"],\n  [\"989\u00f77=141, 2\", \"552\u00f78=69, 0\"],\n  [\"840\u00f77=120, 0\", \"310\u00f72=155, 0\"],\n  [\"881\u00f74=220, 1\", \"147\u00f74=36, 3\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n  @(\"2024-01-22 Monday\", \"2024-01-23 Tuesday\"),\n  @(\"604\u00f73=201, 1\", \"908\u00f72=454, 0\"),\n  @(\"961\u00f75=192, 1\", \"884\u00f73=294, 2\"),\n  @(\"951\u00f79=105, 6\", \"790\u00f78=98, 6\"),\n  @(\"944\u00f72=472, 0\", \"620\u00f78=77, 4\"),\n  @(\"761\u00f79=84, 5\", \"963\u00f73=321, 0\"),\n  @(\"976\u00f78=122, 0\", \"117\u00f77=16, 5\"),\n  @(\"309\u00f78=38, 5\", \"729\u00f78=91, 1\"),\n  @(\"610\u00f78=76, 2\", \"775\u00f75=155, 0\"),\n  @(\"209\u00f77=29, 6\", \"435\u00f74=108, 3\"),\n  @(\"831\u00f76=138, 3\", \"933\u00f79=103, 6\"),\n  @(\"561\u00f76=93, 3\", \"166\u00f75=33, 1\"),\n  @(\"298\u00f74=74, 2\", \"431\u00f73=143, 2\"),\n  @(\"359\u00f72=179, 1\", \"409\u00f78=51, 1\"),\n  @(\"901\u00f72=450, 1\", \"830\u00f73=276, 2\"),\n  @(\"962\u00f73=320, 2\", \"289\u00f74=72, 1\"),\n  @(\"142\u00f72=71, 0\", \"443\u00f74=110, 3\"),\n  @(\"373\u00f79=41, 4\", \"633\u00f77=90, 3\"),\n  @(\"770\u00f74=192, 2\", \"751\u00f75=150, 1\"),\n  @(\"326\u00f73=108, 2\", \"325\u00f74=81, 1\"),\n  @(\"734\u00f78=91, 6\", \"140\u00f76=23, 2\"),\n  @(\"228\u00f78=28, 4\", \"171\u00f77=24, 3\"),\n  @(\"810\u00f79=90, 0\", \"488\u00f78=61, 0\"),\n  @(\"989\u00f77=141, 2\", \"552\u00f78=69, 0\"),\n  @(\"840\u00f77=120, 0\", \"310\u00f72=155, 0\"),\n  @(\"881\u00f74=220, 1\", \"147\u00f74=36, 3\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $result = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n  if (-not $result) {\n    throw \"Text not found: $oldText\"\n  }\n}"}
